$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44769
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 20000
$ws.Range("L2").Value = 20000
$ws.Range("M2").Value = 20000
$ws.Range("P2").Value = 1333

# Row 3
$ws.Range("D3").Value = 44525
$ws.Range("J3").Value = 40
$ws.Range("K3").Value = 8000
$ws.Range("L3").Value = 8000
$ws.Range("M3").Value = 8000
$ws.Range("P3").Value = 533

# Row 4
$ws.Range("D4").Value = 44756
$ws.Range("J4").Value = 80
$ws.Range("K4").Value = 20000
$ws.Range("L4").Value = 20000
$ws.Range("M4").Value = 20000
$ws.Range("P4").Value = 1333

# Row 5 unchanged

# Row 6
$ws.Range("D6").Value = 44518
$ws.Range("J6").Value = 50
$ws.Range("K6").Value = 10000
$ws.Range("L6").Value = 10000
$ws.Range("M6").Value = 10000
$ws.Range("P6").Value = 667

# Row 7
$ws.Range("D7").Value = 44757
$ws.Range("J7").Value = 30

# Row 8
$ws.Range("D8").Value = 44767

# Row 9
$ws.Range("D9").Value = 44776
$ws.Range("J9").Value = 80
$ws.Range("K9").Value = 20000
$ws.Range("L9").Value = 20000
$ws.Range("M9").Value = 20000
$ws.Range("P9").Value = 1333

# Row 10
$ws.Range("D10").Value = 44508
$ws.Range("J10").Value = 40
$ws.Range("K10").Value = 10000
$ws.Range("L10").Value = 10000
$ws.Range("M10").Value = 10000
$ws.Range("P10").Value = 667

# Row 11
$ws.Range("D11").Value = 44749

# Row 12
$ws.Range("D12").Value = 44771
$ws.Range("J12").Value = 40
